# Update profit/cost figures across the Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled profit-recalculation runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2039.7
$ws.Range("J17").Value = 2039.7
$ws.Range("L17").Value = 6119.1
$ws.Range("N17").Value = -6455.1
$ws.Range("H19").Value = 1994.1666
$ws.Range("I19").Value = 1404.5
$ws.Range("J19").Value = 3173.5
$ws.Range("K19").Value = 1404.5
$ws.Range("L19").Value = 3173.5
$ws.Range("M19").Value = -1229.5
$ws.Range("N19").Value = -3523.5
$ws.Range("I28").Value = 859.1429000000001
$ws.Range("J28").Value = 2995
$ws.Range("K28").Value = 859.1429000000001
$ws.Range("L28").Value = 2995
$ws.Range("M28").Value = -374.1429000000001
$ws.Range("N28").Value = -3965
$ws.Range("H33").Value = 116.1
$ws.Range("J33").Value = 272
$ws.Range("L33").Value = 272
$ws.Range("N33").Value = -730
$ws.Range("H41").Value = 322.25
$ws.Range("I41").Value = 322.25
$ws.Range("K41").Value = 322.25
$ws.Range("M41").Value = 117.75
$ws.Range("J62").Value = 9000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10248
$ws.Range("J65").Value = 9000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51240
$ws.Range("H98").Value = 801.7143
$ws.Range("I98").Value = 322.4
$ws.Range("K98").Value = 322.4
$ws.Range("M98").Value = 1175.6
$ws.Range("H122").Value = 801.7143
$ws.Range("I122").Value = 322.4
$ws.Range("K122").Value = 967.1999999999999
$ws.Range("M122").Value = 1482.8
$ws.Range("H125").Value = 4418.3
$ws.Range("I125").Value = 4030.5
$ws.Range("K125").Value = 36274.5
$ws.Range("M125").Value = -33814.5
$ws.Range("H131").Value = 10000
$ws.Range("I131").Value = 10000
$ws.Range("K131").Value = 30000
$ws.Range("M131").Value = -24960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4734.846
$ws.Range("I32").Value = 4296
$ws.Range("K32").Value = 4296
$ws.Range("M32").Value = -4009
$ws.Range("H45").Value = 2238.3333
$ws.Range("I45").Value = 1651.4445
$ws.Range("K45").Value = 1651.4445
$ws.Range("M45").Value = -1274.4445
$ws.Range("H61").Value = 1919
$ws.Range("I61").Value = 785
$ws.Range("K61").Value = 785
$ws.Range("M61").Value = -573
$ws.Range("H74").Value = 11108595
$ws.Range("I74").Value = 13329114
$ws.Range("K74").Value = 13329114
$ws.Range("M74").Value = -13328240
$ws.Range("H77").Value = 11108595
$ws.Range("I77").Value = 13329114
$ws.Range("K77").Value = 66645570
$ws.Range("M77").Value = -66641202
$ws.Range("H102").Value = 1576.3334
$ws.Range("I102").Value = 1514.75
$ws.Range("J102").Value = 1699.5
$ws.Range("K102").Value = 1514.75
$ws.Range("L102").Value = 1699.5
$ws.Range("M102").Value = 107.25
$ws.Range("N102").Value = -4943.5
$ws.Range("H110").Value = 1565.8
$ws.Range("I110").Value = 1565.8
$ws.Range("K110").Value = 1565.8
$ws.Range("M110").Value = 479.2
$ws.Range("H132").Value = 2469.25
$ws.Range("I132").Value = 1544.4
$ws.Range("J132").Value = 3625.3125
$ws.Range("K132").Value = 4633.200000000001
$ws.Range("L132").Value = 10875.9375
$ws.Range("M132").Value = -2103.200000000001
$ws.Range("N132").Value = -15935.9375
$ws.Range("H136").Value = 1919
$ws.Range("I136").Value = 785
$ws.Range("K136").Value = 2355
$ws.Range("M136").Value = 195

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2499.5
$ws.Range("I20").Value = 2499.5
$ws.Range("K20").Value = 2499.5
$ws.Range("M20").Value = -2252.5
$ws.Range("H86").Value = 4507.3335
$ws.Range("J86").Value = 5340.7144
$ws.Range("L86").Value = 5340.7144
$ws.Range("N86").Value = -7586.7144
$ws.Range("H89").Value = 4507.3335
$ws.Range("J89").Value = 5340.7144
$ws.Range("L89").Value = 26703.572
$ws.Range("N89").Value = -37935.572
$ws.Range("H131").Value = 30000
$ws.Range("J131").Value = 30000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("H134").Value = 2115.9092
$ws.Range("I134").Value = 1927.5
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5782.5
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3247.5
$ws.Range("N134").Value = -17070
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H14").Value = 11000
$ws.Range("J14").Value = 11000
$ws.Range("L14").Value = 11000
$ws.Range("N14").Value = -11340
$ws.Range("H58").Value = 2665.125
$ws.Range("I58").Value = 2765
$ws.Range("J58").Value = 2605.2
$ws.Range("K58").Value = 2765
$ws.Range("L58").Value = 2605.2
$ws.Range("M58").Value = -2562
$ws.Range("N58").Value = -3011.2
$ws.Range("H68").Value = 34118
$ws.Range("J68").Value = 34118
$ws.Range("L68").Value = 34118
$ws.Range("N68").Value = -35616
$ws.Range("H71").Value = 34118
$ws.Range("J71").Value = 34118
$ws.Range("L71").Value = 102354
$ws.Range("N71").Value = -109842
$ws.Range("H134").Value = 2972.1667
$ws.Range("I134").Value = 2972.1667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8916.500100000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6381.500100000001
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2665.125
$ws.Range("I136").Value = 2765
$ws.Range("J136").Value = 2605.2
$ws.Range("K136").Value = 8295
$ws.Range("L136").Value = 7815.599999999999
$ws.Range("M136").Value = -5745
$ws.Range("N136").Value = -12915.6
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H92").Value = 709.8
$ws.Range("I92").Value = 599.6667
$ws.Range("J92").Value = 875
$ws.Range("K92").Value = 1799.0001
$ws.Range("L92").Value = 2625
$ws.Range("M92").Value = -551.0001
$ws.Range("N92").Value = -5121
$ws.Range("H112").Value = 7169.8
$ws.Range("I112").Value = 896
$ws.Range("K112").Value = 2688
$ws.Range("M112").Value = -1580
$ws.Range("H131").Value = 1384.16
$ws.Range("I131").Value = 982
$ws.Range("J131").Value = 1511.1578
$ws.Range("K131").Value = 2946
$ws.Range("L131").Value = 4533.4734
$ws.Range("M131").Value = 2094
$ws.Range("N131").Value = -14613.4734
$ws.Range("H132").Value = 3553.111
$ws.Range("I132").Value = 2830
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 25470
$ws.Range("L132").Value = 44994.0015
$ws.Range("M132").Value = -22940
$ws.Range("N132").Value = -50054.0015
$ws.Range("H140").Value = 2599.6667
$ws.Range("I140").Value = 2599.6667
$ws.Range("K140").Value = 7799.000100000001
$ws.Range("M140").Value = -2619.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3799
$ws.Range("I80").Value = 2737.4
$ws.Range("J80").Value = 6453
$ws.Range("K80").Value = 2737.4
$ws.Range("L80").Value = 6453
$ws.Range("M80").Value = -1739.4
$ws.Range("N80").Value = -8449
$ws.Range("H83").Value = 3799
$ws.Range("I83").Value = 2737.4
$ws.Range("J83").Value = 6453
$ws.Range("K83").Value = 13687
$ws.Range("L83").Value = 32265
$ws.Range("M83").Value = -8695
$ws.Range("N83").Value = -42249
$ws.Range("H122").Value = 1685.8572
$ws.Range("I122").Value = 1478.8
$ws.Range("K122").Value = 4436.4
$ws.Range("M122").Value = -1986.4
$ws.Range("H126").Value = 690
$ws.Range("I126").Value = 690
$ws.Range("J126").Value = 690
$ws.Range("K126").Value = 2070
$ws.Range("L126").Value = 2070
$ws.Range("M126").Value = 400
$ws.Range("N126").Value = -7010

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4533
$ws.Range("I7").Value = 4299.5
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 4299.5
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4187.5
$ws.Range("N7").Value = -5224
$ws.Range("H40").Value = 3253.889
$ws.Range("I40").Value = 3416.875
$ws.Range("K40").Value = 3416.875
$ws.Range("M40").Value = -3280.875
$ws.Range("H61").Value = 4978
$ws.Range("I61").Value = 4978
$ws.Range("K61").Value = 4978
$ws.Range("M61").Value = -4776
$ws.Range("H113").Value = 4978
$ws.Range("I113").Value = 4978
$ws.Range("K113").Value = 4978
$ws.Range("M113").Value = -2808
$ws.Range("H126").Value = 4533
$ws.Range("I126").Value = 4299.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 12898.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -10428.5
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 3481.2222
$ws.Range("I132").Value = 2697.0908
$ws.Range("K132").Value = 8091.2724
$ws.Range("M132").Value = -5561.2724
$ws.Range("H136").Value = 13331328
$ws.Range("I136").Value = 13331328
$ws.Range("K136").Value = 39993984
$ws.Range("M136").Value = -39991434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1339.75
$ws.Range("I126").Value = 1295.7333
$ws.Range("K126").Value = 3887.199900000001
$ws.Range("M126").Value = -1417.199900000001
